# This script applies a rearrangement of data rows 3-7 on the active sheet.
# The underlying records (identified by column A id) got reshuffled into
# different rows while keeping columns T:AY (which are identical across
# these rows) untouched. Only columns A,B,D,E,F,G,H,M,Q,R actually change
# value as a net effect of the reshuffle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (target) values for rows 3..7, columns A,B,D,E,F,G,H,M,Q,R
$data = @{
    3 = @{ A=111742294; B=94134; D='NT'; E=53;     F='Vedtrappmossa';        G='Crossocalyx hellerianus'; H='(Nees ex Lindenb.) Meyl.'; M=$null;               Q=331799.9927276275; R=6626510.806996167 }
    4 = @{ A=111742281; B=4711;  D='LC'; E=100299;  F='Thomsons trägnagare';  G='Cacotemnus thomsoni';     H='(Kraatz, 1881)';           M='färska gnagspår';   Q=331821.5503750234; R=6626517.909892835 }
    5 = @{ A=111742269; B=94134; D='NT'; E=53;     F='Vedtrappmossa';        G='Crossocalyx hellerianus'; H='(Nees ex Lindenb.) Meyl.'; M=$null;               Q=331779.9179887357; R=6626525.342625097 }
    6 = @{ A=111742278; B=94134; D='NT'; E=53;     F='Vedtrappmossa';        G='Crossocalyx hellerianus'; H='(Nees ex Lindenb.) Meyl.'; M=$null;               Q=331818.8411813352; R=6626525.099085328 }
    7 = @{ A=111742299; B=94134; D='NT'; E=53;     F='Vedtrappmossa';        G='Crossocalyx hellerianus'; H='(Nees ex Lindenb.) Meyl.'; M=$null;               Q=331807.7707727421; R=6626503.893626045 }
}

foreach ($row in 3..7) {
    $rec = $data[$row]

    $ws.Range("A$row").Value = $rec.A
    $ws.Range("B$row").Value = $rec.B
    $ws.Range("D$row").Value = $rec.D
    $ws.Range("E$row").Value = $rec.E
    $ws.Range("F$row").Value = $rec.F
    $ws.Range("G$row").Value = $rec.G
    $ws.Range("H$row").Value = $rec.H

    if ($rec.M) {
        $ws.Range("M$row").Value = $rec.M
    } else {
        $ws.Range("M$row").Value = $null
    }

    $ws.Range("Q$row").Value = $rec.Q
    $ws.Range("R$row").Value = $rec.R
}
